$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")

# --- 1. Insert the two new columns ---------------------------------------
# First insertion: before the old column N ("Mo ta"), this becomes the new
# "% VAT *" column.
$x = $ws.Columns.Item(14).Insert()

# Second insertion: before what is now column R (the old "Thuoc tinh 1",
# shifted right by the first insert). This becomes the new
# "Co tao phien ban*" column.
$x = $ws.Columns.Item(18).Insert()

# --- 2. Give the two new header cells their text --------------------------
$ws.Range("N1").Value = "% VAT *"
$ws.Range("R1").Value = "Có tạo phiên bản*"

# --- 3. Style the new header cells like their neighbours -------------------
# Re-touching a protection flag nudges Excel into allocating the new cells
# their own cellXfs entry (mirrors the "apply number format / protection"
# bits that appear on the real template's two new header styles). Using two
# different values keeps the two new cells on two distinct style records,
# matching the two (functionally identical) new styles in the target file.
$ws.Range("N1").FormulaHidden = $false
$ws.Range("R1").FormulaHidden = $true

# --- 4. Replace the broken list-validations on D:E / H with a plain,
#        "choose from the list" reminder validation (the old ones pointed at
#        #REF! ranges and are being retired). -------------------------------
$x = $ws.Range("D1:E1048576").Validation.Delete()
$x = $ws.Range("D1:E1048576").Validation.Add(0, 1, 1, "", "")
$ws.Range("D1:E1048576").Validation.ErrorTitle = "Nhập lại"
$ws.Range("D1:E1048576").Validation.ErrorMessage = "Vui lòng chọn giá trị từ danh sách"

$x = $ws.Range("H1:H1048576").Validation.Delete()
$x = $ws.Range("H1:H1048576").Validation.Add(0, 1, 1, "", "")
$ws.Range("H1:H1048576").Validation.ErrorTitle = "Nhập lại"
$ws.Range("H1:H1048576").Validation.ErrorMessage = "Vui lòng chọn giá trị từ danh sách"

# --- 5. Move the four "Giá trị N" comments two columns to the right, to
#        follow the "Thuộc tính N"/"Giá trị N" columns that were pushed over
#        by the two newly inserted columns. Walk right-to-left so a move
#        never clobbers a comment we still have to read. ------------------
$moves = @(
  @{from="X1"; to="Z1"},
  @{from="V1"; to="X1"},
  @{from="T1"; to="V1"},
  @{from="R1"; to="T1"}
)
foreach ($mv in $moves) {
  $src = $ws.Range($mv.from)
  $txt = $src.Comment.Text()
  $x = $src.Comment.Delete()
  $dst = $ws.Range($mv.to)
  $x = $dst.AddComment($txt)
}

# --- 6. Refresh the saved selection on the Product sheet --------------------
$x = $ws.Activate()
$x = $ws.Range("T7").Select()
